# Insert a new weekly price record for "Apio" (Vega Modelo de Temuco) above the
# existing row 134, shifting the subsequent rows (134-179) down to (135-180).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 134 (Excel default shift = down), which also pushes
# the used range/dimension from R179 to R180 automatically.
$ws.Rows.Item(134).Insert()

# Populate the newly inserted row 134 with this week's data.
$ws.Cells.Item(134, 1).Value  = 10
$ws.Cells.Item(134, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(134, 3).Value  = "La Araucanía"
$ws.Cells.Item(134, 4).Value  = 44468
$ws.Cells.Item(134, 5).Value  = 9
$ws.Cells.Item(134, 6).Value  = 100112017
$ws.Cells.Item(134, 7).Value  = "Apio"
$ws.Cells.Item(134, 8).Value  = "Americana (o)"
$ws.Cells.Item(134, 9).Value  = "Primera"
$ws.Cells.Item(134, 10).Value = 30
$ws.Cells.Item(134, 11).Value = 10000
$ws.Cells.Item(134, 12).Value = 10000
$ws.Cells.Item(134, 13).Value = 10000
$ws.Cells.Item(134, 14).Value = "$/docena de matas"
$ws.Cells.Item(134, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(134, 16).Value = 1667
$ws.Cells.Item(134, 17).Value = 6
$ws.Cells.Item(134, 18).Value = "Hortaliza"
